$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Enterprises density (per 1000 people)" row (row 13): update Micro/SMEs/MSMEs
# values with more precise decimals. These must stay text (as in the source
# workbook) rather than being auto-coerced to numbers by Excel, so we mark
# the cells as Text-formatted before writing the new values.
$rng = $ws.Range("B13:D13")
$rng.NumberFormat = "@"

$ws.Range("B13").Value = "12.12"
$ws.Range("C13").Value = "3.93"
$ws.Range("D13").Value = "16.04"
